$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1729.6086
$ws.Range("J112").Value = 1979.4706
$ws.Range("L112").Value = 5938.4118
$ws.Range("N112").Value = -8154.4118
$ws.Range("H132").Value = 9117.200000000001
$ws.Range("I132").Value = 9035.5
$ws.Range("K132").Value = 27106.5
$ws.Range("M132").Value = -24576.5
$ws.Range("H137").Value = 6899.1924
$ws.Range("I137").Value = 3669.3684
$ws.Range("J137").Value = 15665.857
$ws.Range("K137").Value = 11008.1052
$ws.Range("L137").Value = 46997.571
$ws.Range("M137").Value = -8458.1052
$ws.Range("N137").Value = -52097.571
$ws.Range("H138").Value = 2584.9744
$ws.Range("J138").Value = 2785.6033
$ws.Range("L138").Value = 8356.8099
$ws.Range("N138").Value = -18636.8099
$ws.Range("H141").Value = 13953
$ws.Range("I141").Value = 13953
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 41859
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -36679
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 18299316
$ws.Range("I61").Value = 13894486
$ws.Range("K61").Value = 13894486
$ws.Range("M61").Value = -13894274
$ws.Range("H136").Value = 18299316
$ws.Range("I136").Value = 13894486
$ws.Range("K136").Value = 41683458
$ws.Range("M136").Value = -41680908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 313743.12
$ws.Range("I134").Value = 1123.9642
$ws.Range("J134").Value = 2502077.2
$ws.Range("K134").Value = 3371.8926
$ws.Range("L134").Value = 7506231.600000001
$ws.Range("M134").Value = -836.8925999999997
$ws.Range("N134").Value = -7511301.600000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4082.9473
$ws.Range("I94").Value = 6001.6665
$ws.Range("J94").Value = 3723.1875
$ws.Range("K94").Value = 6001.6665
$ws.Range("L94").Value = 3723.1875
$ws.Range("M94").Value = -5550.6665
$ws.Range("N94").Value = -4625.1875
$ws.Range("H99").Value = 1736.0588
$ws.Range("I99").Value = 1575
$ws.Range("J99").Value = 2122.6
$ws.Range("K99").Value = 1575
$ws.Range("L99").Value = 2122.6
$ws.Range("M99").Value = -77
$ws.Range("N99").Value = -5118.6
$ws.Range("H122").Value = 963.3333
$ws.Range("I122").Value = 1061.3334
$ws.Range("J122").Value = 473.33334
$ws.Range("K122").Value = 3184.0002
$ws.Range("L122").Value = 1420.00002
$ws.Range("M122").Value = -734.0001999999999
$ws.Range("N122").Value = -6320.000019999999
$ws.Range("H126").Value = 1736.0588
$ws.Range("I126").Value = 1575
$ws.Range("J126").Value = 2122.6
$ws.Range("K126").Value = 4725
$ws.Range("L126").Value = 6367.799999999999
$ws.Range("M126").Value = -2255
$ws.Range("N126").Value = -11307.8
$ws.Range("H134").Value = 2909.84
$ws.Range("I134").Value = 1319.7778
$ws.Range("J134").Value = 6998.5713
$ws.Range("K134").Value = 3959.3334
$ws.Range("L134").Value = 20995.7139
$ws.Range("M134").Value = -1424.3334
$ws.Range("N134").Value = -26065.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1899.4117
$ws.Range("I5").Value = 1899.4117
$ws.Range("K5").Value = 5698.2351
$ws.Range("M5").Value = -5586.2351
$ws.Range("H12").Value = 1185.5454
$ws.Range("J12").Value = 299.17648
$ws.Range("L12").Value = 897.52944
$ws.Range("N12").Value = -1243.52944
$ws.Range("H68").Value = 4998.25
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4998.25
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14994.75
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -16616.75
$ws.Range("H69").Value = 3308.3333
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 3962.5
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 11887.5
$ws.Range("M69").Value = -5189
$ws.Range("N69").Value = -13509.5
$ws.Range("H70").Value = 142858510
$ws.Range("I70").Value = 142858510
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 428575530
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -428575215
$ws.Range("N70").ClearContents()
$ws.Range("H71").Value = 4998.25
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4998.25
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 44984.25
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -53096.25
$ws.Range("H72").Value = 3308.3333
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 3962.5
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 35662.5
$ws.Range("M72").Value = -13944
$ws.Range("N72").Value = -43774.5
$ws.Range("H73").Value = 142858510
$ws.Range("I73").Value = 142858510
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 428575530
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -428574438
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 4050.0625
$ws.Range("J80").Value = 4093.2666
$ws.Range("L80").Value = 12279.7998
$ws.Range("N80").Value = -14151.7998
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 4050.0625
$ws.Range("J83").Value = 4093.2666
$ws.Range("L83").Value = 36839.3994
$ws.Range("N83").Value = -46199.3994
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H92").Value = 301
$ws.Range("I92").Value = 302
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 906
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 342
$ws.Range("N92").Value = -3396
$ws.Range("H122").Value = 2601
$ws.Range("I122").Value = 523.55554
$ws.Range("J122").Value = 3700.8235
$ws.Range("K122").Value = 4711.99986
$ws.Range("L122").Value = 33307.4115
$ws.Range("M122").Value = -2261.99986
$ws.Range("N122").Value = -38207.4115
$ws.Range("H131").Value = 3363.8906
$ws.Range("I131").Value = 570.8570999999999
$ws.Range("J131").Value = 4145.94
$ws.Range("K131").Value = 1712.5713
$ws.Range("L131").Value = 12437.82
$ws.Range("M131").Value = 3327.4287
$ws.Range("N131").Value = -22517.82
$ws.Range("H132").Value = 1743.5
$ws.Range("I132").Value = 854
$ws.Range("J132").Value = 1891.75
$ws.Range("K132").Value = 7686
$ws.Range("L132").Value = 17025.75
$ws.Range("M132").Value = -5156
$ws.Range("N132").Value = -22085.75
$ws.Range("H135").Value = 1899.4117
$ws.Range("I135").Value = 1899.4117
$ws.Range("K135").Value = 17094.7053
$ws.Range("M135").Value = -14559.7053
$ws.Range("H137").Value = 6042.778
$ws.Range("I137").Value = 3249.75
$ws.Range("J137").Value = 8277.200000000001
$ws.Range("K137").Value = 9749.25
$ws.Range("L137").Value = 24831.6
$ws.Range("M137").Value = -4649.25
$ws.Range("N137").Value = -35031.60000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4572.375
$ws.Range("I126").Value = 789.5
$ws.Range("K126").Value = 2368.5
$ws.Range("M126").Value = 101.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38650.766
$ws.Range("I7").Value = 4629.095
$ws.Range("J7").Value = 118034.664
$ws.Range("K7").Value = 4629.095
$ws.Range("L7").Value = 118034.664
$ws.Range("M7").Value = -4517.095
$ws.Range("N7").Value = -118258.664
$ws.Range("H9").Value = 1505.5
$ws.Range("J9").Value = 1505.5
$ws.Range("L9").Value = 1505.5
$ws.Range("N9").Value = -1953.5
$ws.Range("H12").Value = 67627
$ws.Range("J12").Value = 90002.664
$ws.Range("L12").Value = 90002.664
$ws.Range("N12").Value = -90342.664
$ws.Range("H126").Value = 38650.766
$ws.Range("I126").Value = 4629.095
$ws.Range("J126").Value = 118034.664
$ws.Range("K126").Value = 13887.285
$ws.Range("L126").Value = 354103.992
$ws.Range("M126").Value = -11417.285
$ws.Range("N126").Value = -359043.992
